$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab20")

# Fix mojibake/encoding issue in the footnote about Regional Economic Communities (cell A103)
$ws.Range("A103").Value = "Regional Economic Communities:CEN-SAD = ""Community of Sahel-Saharan States"";COMESA = ""Common Market for Eastern and Southern Africa"";EAC = ""East African Community"";ECCAS = ""Economic Community of Central African States"";ECOWAS = ""Economic Community of West African States"";IGAD = ""Intergovernmental Authority on Development"";SADC = ""Southern African Development Community"";UMA = ""Arab Maghreb Union"";PALOP = ""Países Africanos de Língua Oficial Portuguesa"";ASEAN = ""Association of Southeast Asian Nations"";MERCOSUR = ""Mercado Común del Sur"".EU27 = ""European Union (27 members)"".OECD = ""Organisation for Economic Co-operation and Development""."

# Update data values for row 97 (Africa, Fragile States)
$ws.Range("C97").Value = 1.5251825098374601
$ws.Range("D97").Value = 1.1256550052454499
$ws.Range("E97").Value = 2.76671116250446
$ws.Range("F97").Value = 2.8827277667316
$ws.Range("G97").Value = 8.30027644431898
$ws.Range("H97").Value = 1.3712836397870001
$ws.Range("I97").Value = 3.0293076778924499
$ws.Range("J97").Value = 17791.2819892468
$ws.Range("K97").Value = 8951.2158386914307
$ws.Range("L97").Value = 33477.35
$ws.Range("M97").Value = 35073.7514204974
$ws.Range("N97").Value = 95293.599248435698
$ws.Range("O97").Value = 12711.666735458801
$ws.Range("P97").Value = 28120.1522129474

# Update data values for row 98 (ROW, Fragile States)
$ws.Range("C98").Value = 1.16286921245604
$ws.Range("D98").Value = 0.44740303145552002
$ws.Range("E98").Value = 0.96299166631921995
$ws.Range("F98").Value = 7.7596505474386603
$ws.Range("G98").Value = 10.3329144576694
$ws.Range("H98").Value = 1.31621981515466
$ws.Range("I98").Value = 3.9879152073316901
$ws.Range("J98").Value = 19452.905651664201
$ws.Range("K98").Value = 4312.2249860585498
$ws.Range("L98").Value = 15907.48
$ws.Range("M98").Value = 84787.025146911605
$ws.Range("N98").Value = 124459.635784634
$ws.Range("O98").Value = 13920.4952674058
$ws.Range("P98").Value = 37587.360629103401
